$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4198.545839504756
$ws.Range("C3").Value = 4198.545839504756
$ws.Range("C4").Value = 4159.108552782584
$ws.Range("C5").Value = 4116.962347901228
$ws.Range("C6").Value = 3951.105791765196
$ws.Range("C7").Value = 3951.105791765196
$ws.Range("C8").Value = 3889.409402516725
$ws.Range("C9").Value = 3889.409402516725
$ws.Range("C10").Value = 3889.409402516725
$ws.Range("C11").Value = 3819.228336323561
$ws.Range("C12").Value = 3767.663951912773
